$wb = $excel.ActiveWorkbook

# --- Sheet 1: Liquidity Gap ---
$ws1 = $wb.Worksheets.Item("Liquidity Gap")
$ws1.Range("B5").Value = 40097834
$ws1.Range("C5").Value = 17650714
$ws1.Range("D5").Value = 22447120
$ws1.Range("E5").Value = 35539955
$ws1.Range("B6").Value = 50865822
$ws1.Range("C6").Value = 912681
$ws1.Range("D6").Value = 49953141
$ws1.Range("E6").Value = 85493096
$ws1.Range("B7").Value = 181047468
$ws1.Range("C7").Value = 22049274
$ws1.Range("D7").Value = 158998194
$ws1.Range("E7").Value = 244491290
$ws1.Range("B8").Value = 137703222
$ws1.Range("C8").Value = 95188050
$ws1.Range("D8").Value = 42515172
$ws1.Range("E8").Value = 287006462
$ws1.Range("B9").Value = 48148807
$ws1.Range("C9").Value = 67867558
$ws1.Range("D9").Value = -19718751
$ws1.Range("E9").Value = 267287711
$ws1.Range("B10").Value = 0
$ws1.Range("C10").Value = 96621244
$ws1.Range("D10").Value = -96621244
$ws1.Range("E10").Value = 170666467

# --- Sheet 2: Repricing Gap ---
$ws2 = $wb.Worksheets.Item("Repricing Gap")
$ws2.Range("B4").Value = 48871432
$ws2.Range("C4").Value = 78670
$ws2.Range("D4").Value = 48792762
$ws2.Range("E4").Value = 55189425
$ws2.Range("B5").Value = 90462007
$ws2.Range("C5").Value = 33520374
$ws2.Range("D5").Value = 56941633
$ws2.Range("E5").Value = 112131058
$ws2.Range("B6").Value = 67725781
$ws2.Range("C6").Value = 31488965
$ws2.Range("D6").Value = 36236816
$ws2.Range("E6").Value = 148367874
$ws2.Range("B7").Value = 128817688
$ws2.Range("C7").Value = 46327208
$ws2.Range("D7").Value = 82490480
$ws2.Range("E7").Value = 230858354
$ws2.Range("B8").Value = 99349448
$ws2.Range("C8").Value = 101287781
$ws2.Range("D8").Value = -1938333
$ws2.Range("E8").Value = 228920021
$ws2.Range("B9").Value = 29123514
$ws2.Range("C9").Value = 37233029
$ws2.Range("D9").Value = -8109515
$ws2.Range("E9").Value = 220810506
$ws2.Range("B10").Value = 0
$ws2.Range("C10").Value = 50144039
$ws2.Range("D10").Value = -50144039
$ws2.Range("E10").Value = 170666467

# --- Sheet 3: NII Sensitivity ---
$ws3 = $wb.Worksheets.Item("NII Sensitivity")
$ws3.Range("A1").Clear()
$ws3.Range("B3").Value = 2308583.54
$ws3.Range("B4").Value = -2308583.54
$ws3.Range("B5").Value = 230858354
